$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=9;  A="SP22092022001215"; B=100; C="PANADOL STRIP 10"; D=1; E=15 },
    @{ Row=10; A="SP22092022001215"; B=101; C="PANADOL STRIP 20"; D=1; E=29 },
    @{ Row=11; A="SP22092022001215"; B=102; C="NUROFEN STRIP 15"; D=1; E=12 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
